# Individual Contribution Metrics - "Added second delivery files"
# Fill in the E/F/G (per-member score) columns for the rows that were
# missing them, add a new underlined/empty marker cell at D13, move the
# active selection, and refresh the page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing per-member score columns (E, F, G) -------------------
# Row 6: only F/G were filled in; add the missing E value.
$ws.Range("E6").Value = 5

# Row 7: only C was filled in; add E/F/G to match C.
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 5

# Row 9: only C was filled in; add E/F/G to match C.
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 13

# Row 10: only C was filled in; add E/F/G to match C.
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5

# --- New underlined style used as a marker in D13 --------------------------
$ws.Range("D13").Font.Underline = $true

# --- Update the selection to match the author's last-saved cursor ---------
$ws.Range("D13").Select()

# --- Refresh page setup (adds printer settings / portrait orientation) ----
$ws.PageSetup.Orientation = 1
